$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the SONG_ID value for row 2 (shared string formerly "clscam" -> "batimt")
$ws.Range("A2").Value = "batimt"

# Clear the other data values on row 2 (GENRE, ONI/NORMAL/HARD/EDIT/EASY level columns)
$ws.Range("B2").Clear()
$ws.Range("E2:I2").Clear()

# Remove rows 3 and 4 entirely (they held the "crtesc" / "fmod" songs)
$ws.Rows("3:4").Delete()

# Update selection to match the saved state
$ws.Range("A2").Select()
